# Daily attendance processing - 2025-10-10 22:20:02
# Reorders the "Recorded By" (column G) email lists for the affected rows
# so that the join order matches the freshly generated report.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 'mennatulla.medhat@med.asu.edu.eg, servinaz@med.asu.edu.eg'
$ws.Range("G3").Value = 'nourhan.mahmoud@med.asu.edu.eg , servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg'
$ws.Range("G9").Value = 'wessam.atef@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg'
$ws.Range("G15").Value = 'Mayarembaby@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg, nourhanmamdouh@med.asu.edu.eg'
$ws.Range("G27").Value = 'System, nourhanmohamed@med.asu.edu.eg, marian.samir@med.asu.edu.eg, randa.rabea@med.asu.edu.eg'
$ws.Range("G29").Value = 'yasmin.m.senosy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marina_atef@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, marinasorial@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg'
$ws.Range("G33").Value = 'mennatulla.medhat@med.asu.edu.eg, servinaz@med.asu.edu.eg'
$ws.Range("G34").Value = 'nourhan.mahmoud@med.asu.edu.eg , servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg'
$ws.Range("G40").Value = 'wessam.atef@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg'
$ws.Range("G46").Value = 'afnan.fares@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg'
$ws.Range("G58").Value = 'System, nourhanmohamed@med.asu.edu.eg, marian.samir@med.asu.edu.eg, randa.rabea@med.asu.edu.eg'
$ws.Range("G60").Value = 'yasmin.m.senosy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marina_atef@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, marinasorial@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg'
$ws.Range("G64").Value = 'mennatulla.medhat@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg'
$ws.Range("G65").Value = 'mennatulla.medhat@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Range("G71").Value = 'mariam.noureldin@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg'
$ws.Range("G76").Value = 'Mayarembaby@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg, nourhanmamdouh@med.asu.edu.eg'
$ws.Range("G77").Value = 'afnan.fares@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg'
$ws.Range("G89").Value = 'nourhan.mostafa@med.asu.edu.eg, Aya.alaa-eldein@med.asu.edu.eg'
$ws.Range("G91").Value = 'yasmin.m.senosy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marina_atef@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, marinasorial@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg'
$ws.Range("G95").Value = 'mennatulla.medhat@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg'
$ws.Range("G96").Value = 'alshimaa.atef@med.asu.edu.egm, nourhan.mahmoud@med.asu.edu.eg , servinaz@med.asu.edu.eg'
$ws.Range("G102").Value = 'norhan.mohamed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, aml.awwad@med.asu.edu.eg'
$ws.Range("G108").Value = 'Mayarembaby@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg, nourhanmamdouh@med.asu.edu.eg'
$ws.Range("G122").Value = 'yasmin.m.senosy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg'
$ws.Range("G126").Value = 'mennatulla.medhat@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg'
$ws.Range("G127").Value = 'alshimaa.atef@med.asu.edu.egm, nourhan.mahmoud@med.asu.edu.eg , servinaz@med.asu.edu.eg'
$ws.Range("G133").Value = 'norhan.mohamed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, aml.awwad@med.asu.edu.eg'
$ws.Range("G138").Value = 'Mayarembaby@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg, nourhanmamdouh@med.asu.edu.eg'
$ws.Range("G139").Value = 'Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg'
$ws.Range("G140").Value = 'afnan.fares@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg'
$ws.Range("G151").Value = "Amr-Saeed@med.asu.edu.eg`n, enas.omran@med.asu.edu.eg, takyas@med.asu.edu.eg, user@user.com"
$ws.Range("G153").Value = 'yasmin.m.senosy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg'
$ws.Range("G157").Value = 'eman.tantawi@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg'
$ws.Range("G158").Value = 'mennatulla.medhat@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Range("G164").Value = 'mariam.noureldin@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg'
$ws.Range("G170").Value = 'Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg'
$ws.Range("G184").Value = 'wafaa.ebida@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg'
